$d = $word.ActiveDocument

$newValues = @("17+9=","9+12=","50-14=","86+5=","34+38=","67+28=","72+9=","81-62=","37+57=","54+9=","91-7=","4+67=","80-15=","85-49=","9+8=","18+8=","25+39=","11-8=","70-48=","18+48=","66-9=","39+19=","36+8=","39+18=","45+8=","6+56=","45+47=","31-12=","82-15=","69+26=","79+17=","94-47=","91-48=","35+36=","73-38=","72-56=","92-17=","19+32=","38+7=","70-4=","18+58=","8+69=","19+39=","41-9=","58+18=","29+64=","64-16=","96-47=","67+27=","84-58=","7+9=","37+46=","91-25=","60-8=","6+5=","7+14=","76-47=","25+39=","69+25=","89+3=","81-68=","7+26=","65-36=","29+5=","62-24=","93-44=","19+46=","34+47=","7+26=","14+59=","55+27=","73-56=","85-19=","55+8=","72-59=","61-46=","72+9=","17+56=","91-35=","85-78=","29+29=","7+7=","3+38=","19+17=","8+49=","93-17=","21-12=","50-33=","7+6=","35+46=","90-24=","32-17=","6+37=","8+57=","32-19=","39+45=","69+17=","25+48=","21-7=","65+19=")

$tbl = $d.Tables.Item(1)
$idx = 0
foreach ($row in $tbl.Rows) {
    foreach ($cell in $row.Cells) {
        if ($idx -lt $newValues.Count) {
            $val = $newValues[$idx]
            $r = $cell.Range
            $r.End = $r.End - 1
            $r.Text = $val
        }
        $idx = $idx + 1
    }
}

Write-Output "Replaced $idx cells (expected $($newValues.Count))"
